$wb = $excel.ActiveWorkbook

# The text "2016-08-29 14:50:54" is shared by two cells (Overview!G3 and
# de-de!H3) that happen to hold the same timestamp; both must be updated
# together since they reference the same shared-string entry.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-29 14:52:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-29 14:52:15"
$wsDeDe.Range("K3").Value = "2016-08-29 14:52:47"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-29 14:51:59"
$wsZhCn.Range("K3").Value = "2016-08-29 14:52:39"
